# Update the cryptos price list (Price / Volume(1h) columns) with the
# latest scraped figures, as produced by the GitHub Actions refresh job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell address plus its new text value. Every D/E cell in
# this sheet is stored as text (prices use "."-grouped digits, percentages
# keep their padding spaces), so every update below is written as a string.
$updates = @(
    @{ Addr = 'D2';  Value = '70.147.49' },
    @{ Addr = 'E2';  Value = '  +0.09%  ' },
    @{ Addr = 'D3';  Value = '3.593.35' },
    @{ Addr = 'E3';  Value = '  -0.18%  ' },
    @{ Addr = 'E4';  Value = '  +0.04%  ' },
    @{ Addr = 'D5';  Value = '579.95' },
    @{ Addr = 'E5';  Value = '  -1.54%  ' },
    @{ Addr = 'D6';  Value = '192.33' },
    @{ Addr = 'E6';  Value = '  +0.75%  ' },
    @{ Addr = 'E7';  Value = '  -2.02%  ' },
    @{ Addr = 'D8';  Value = '3.589.37' },
    @{ Addr = 'E8';  Value = '  -0.04%  ' },
    @{ Addr = 'E9';  Value = '  +0.04%  ' },
    @{ Addr = 'E10'; Value = '  +2.71%  ' },
    @{ Addr = 'D11'; Value = '0.666' },
    @{ Addr = 'E11'; Value = '  +0.53%  ' },
    @{ Addr = 'D12'; Value = '56.12' },
    @{ Addr = 'E12'; Value = '  -3.22%  ' },
    @{ Addr = 'E13'; Value = '  +5.09%  ' },
    @{ Addr = 'D14'; Value = '9.66' },
    @{ Addr = 'E14'; Value = '  -1.34%  ' },
    @{ Addr = 'D15'; Value = '4.167.62' },
    @{ Addr = 'E15'; Value = '  -0.03%  ' },
    @{ Addr = 'D16'; Value = '20.06' },
    @{ Addr = 'E16'; Value = '  +3.45%  ' },
    @{ Addr = 'D17'; Value = '3.584.77' },
    @{ Addr = 'E17'; Value = '  -0.30%  ' },
    @{ Addr = 'D18'; Value = '70.142.06' },
    @{ Addr = 'E18'; Value = '  +0.26%  ' },
    @{ Addr = 'E19'; Value = '  +1.74%  ' },
    @{ Addr = 'E20'; Value = '  +0.27%  ' },
    @{ Addr = 'E21'; Value = '  -0.44%  ' },
    @{ Addr = 'D22'; Value = '478.92' },
    @{ Addr = 'E22'; Value = '  -3.34%  ' },
    @{ Addr = 'D23'; Value = '19.55' },
    @{ Addr = 'E23'; Value = '  +10.65%  ' },
    @{ Addr = 'D24'; Value = '5.06' },
    @{ Addr = 'E24'; Value = '  -6.13%  ' },
    @{ Addr = 'D25'; Value = '95.97' },
    @{ Addr = 'E25'; Value = '  +5.77%  ' },
    @{ Addr = 'E26'; Value = '  -1.75%  ' },
    @{ Addr = 'E27'; Value = '  -3.06%  ' },
    @{ Addr = 'D28'; Value = '11.11' },
    @{ Addr = 'E28'; Value = '  -0.20%  ' },
    @{ Addr = 'D29'; Value = '9.47' },
    @{ Addr = 'E29'; Value = '  +0.54%  ' },
    @{ Addr = 'D30'; Value = '32.42' },
    @{ Addr = 'E30'; Value = '  +0.38%  ' },
    @{ Addr = 'D31'; Value = '7.69' },
    @{ Addr = 'E31'; Value = '  +0.38%  ' },
    @{ Addr = 'E32'; Value = '  +0.11%  ' },
    @{ Addr = 'E33'; Value = '  +1.69%  ' },
    @{ Addr = 'D34'; Value = '66.46' },
    @{ Addr = 'E34'; Value = '  +2.03%  ' },
    @{ Addr = 'D35'; Value = '586.91' },
    @{ Addr = 'E35'; Value = '  -4.98%  ' },
    @{ Addr = 'D36'; Value = '39.15' },
    @{ Addr = 'E36'; Value = '  +2.36%  ' },
    @{ Addr = 'E37'; Value = '  +0.01%  ' },
    @{ Addr = 'E38'; Value = '  -1.65%  ' },
    @{ Addr = 'D39'; Value = '0.398' },
    @{ Addr = 'E39'; Value = '  -1.88%  ' },
    @{ Addr = 'D40'; Value = '3.24' },
    @{ Addr = 'E40'; Value = '  +19.80%  ' },
    @{ Addr = 'E41'; Value = '  -5.84%  ' },
    @{ Addr = 'E42'; Value = '  -5.10%  ' },
    @{ Addr = 'D43'; Value = '2.88' },
    @{ Addr = 'E43'; Value = '  +7.59%  ' },
    @{ Addr = 'D44'; Value = '3.240.74' },
    @{ Addr = 'E44'; Value = '  -2.46%  ' },
    @{ Addr = 'E45'; Value = '  +0.55%  ' },
    @{ Addr = 'E46'; Value = '  -0.28%  ' },
    @{ Addr = 'D47'; Value = '3.37' },
    @{ Addr = 'E47'; Value = '  +2.50%  ' },
    @{ Addr = 'D48'; Value = '9.47' },
    @{ Addr = 'E48'; Value = '  +4.05%  ' },
    @{ Addr = 'E49'; Value = '  +0.57%  ' },
    @{ Addr = 'D50'; Value = '0.999' },
    @{ Addr = 'E50'; Value = '  +0.11%  ' },
    @{ Addr = 'D51'; Value = '3.16' },
    @{ Addr = 'E51'; Value = '  -5.06%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    # Column D values such as "20.06" or "0.999" would otherwise be
    # auto-coerced to numbers by the normal Value setter (Excel's usual
    # type-inference on assignment). Force the cell to Text first so the
    # string is stored verbatim, matching the sheet's original inline-text
    # cells, then drop the temporary number format so no stray style is
    # left behind on the cell.
    if ($u.Addr.StartsWith('D')) {
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
